# feat: add 2022-Q3 data
#
# 1. Insert a new worksheet "2022-Q3" right after "总计", containing the
#    fund-holdings breakdown for the new quarter.
# 2. Update the "总计" (summary) sheet: a new top row for 2022-Q3 is
#    inserted and all the quarters below shift down by one row.

$wb = $excel.ActiveWorkbook

# xlPasteSpecial constants
$xlPasteValues  = -4163
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# Step 1: perform ALL worksheet-structure changes first. (Sheet object
# references fetched by name/index before an Add() can go stale/shift,
# so nothing below re-uses a sheet handle obtained before this point.)
# ---------------------------------------------------------------------
$totalSheetTmp = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($null, $totalSheetTmp)
$newSheet.Name = "2022-Q3"

# Now re-fetch every sheet handle we need, fresh, post-restructure.
$totalSheet = $wb.Worksheets.Item("总计")
$q2Sheet    = $wb.Worksheets.Item("2022-Q2")

# Helper: write a value into a cell but force it to stay TEXT even when the
# string looks numeric (e.g. "6.97", "160611", "0.1450"). We stage the text
# in a scratch cell (forced to text with a leading apostrophe), copy just
# the *value* (not the format) onto the destination, then clear the scratch
# cell. This keeps the destination's existing number format/style intact.
function Set-TextValue {
    param($sheet, $cell, [string]$text)
    $scratch = $sheet.Range("ZZ9000")
    $scratch.Value = "'" + $text
    $scratch.Copy()
    $cell.PasteSpecial($xlPasteValues)
    $scratch.Clear()
}

# ---- header row (copy formatting from the 2022-Q2 sheet's header row) ----
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $src = $q2Sheet.Cells.Item(1, $col)
    $dst = $newSheet.Cells.Item(1, $col)
    $dst.Value = $headers[$col - 2]
    $src.Copy()
    $dst.PasteSpecial($xlPasteFormats)
}

# ---- data rows ----
# columns: A=index(n) B=code(text) C=name(text) D=scale(text) E=position(text)
#          F=ratio(text) G=marketvalue(text) H=rank(n)
$rows = @(
    @{ A = 0; B = "160611"; C = "鹏华优质治理混合（LOF）"; D = "6.97"; E = "59.74"; F = "2.08"; G = "0.1450"; H = 10 },
    @{ A = 1; B = "014831"; C = "兴银中证1000指数增强A";    D = "1.37"; E = "83.33"; F = "1.50"; G = "0.0206"; H = 1 },
    @{ A = 2; B = "014832"; C = "兴银中证1000指数增强C";    D = "0.90"; E = "83.33"; F = "1.50"; G = "0.0135"; H = 1 }
)

$r = 2
foreach ($row in $rows) {
    $cellA = $newSheet.Cells.Item($r, 1)
    $srcA  = $q2Sheet.Cells.Item($r, 1)
    $cellA.Value = $row.A
    $srcA.Copy()
    $cellA.PasteSpecial($xlPasteFormats)

    Set-TextValue $newSheet $newSheet.Cells.Item($r, 2) $row.B
    Set-TextValue $newSheet $newSheet.Cells.Item($r, 3) $row.C
    Set-TextValue $newSheet $newSheet.Cells.Item($r, 4) $row.D
    Set-TextValue $newSheet $newSheet.Cells.Item($r, 5) $row.E
    Set-TextValue $newSheet $newSheet.Cells.Item($r, 6) $row.F
    Set-TextValue $newSheet $newSheet.Cells.Item($r, 7) $row.G

    $newSheet.Cells.Item($r, 8).Value = $row.H

    $r++
}

# ---------------------------------------------------------------------
# Step 2: update the "总计" summary sheet - insert the 2022-Q3 row and
# shift the rest of the quarters down
# ---------------------------------------------------------------------
$summaryRows = @(
    @{ B = "2022-Q3"; C = 3;  D = 0.18 },
    @{ B = "2022-Q2"; C = 68; D = 17.18 },
    @{ B = "2022-Q1"; C = 73; D = 37.61 },
    @{ B = "2021-Q4"; C = 82; D = 46.97 },
    @{ B = "2021-Q3"; C = 29; D = 12.83 },
    @{ B = "2021-Q2"; C = 20; D = 6.77 },
    @{ B = "2021-Q1"; C = 4;  D = 0.3 },
    @{ B = "2020-Q4"; C = 3;  D = 0.2 }
)

$r = 2
$idx = 0
foreach ($row in $summaryRows) {
    $cellA = $totalSheet.Cells.Item($r, 1)
    if ($r -gt 8) {
        # row 9 is brand new - copy column-A styling from the row above
        $aboveA = $totalSheet.Cells.Item($r - 1, 1)
        $aboveA.Copy()
        $cellA.PasteSpecial($xlPasteFormats)
    }
    $cellA.Value = $idx

    $totalSheet.Cells.Item($r, 2).Value = $row.B
    $totalSheet.Cells.Item($r, 3).Value = $row.C
    $totalSheet.Cells.Item($r, 4).Value = $row.D

    $r++
    $idx++
}

Write-Host "2022-Q3 sheet added and summary sheet updated"
